# Update data-raw/alt_countries.xlsx:
# Insert a new "altname4" column between "altname3" and "formername",
# and record the historical name "Zaire" / "Congo Democratic Republic"
# for the Democratic Republic of the Congo (COD) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column E (shifts former "formername"/"formername2" columns
# from E/F to F/G).
$ws.Columns("E:E").Insert()

# Give the new column the same custom width as the other alt-name columns.
$ws.Columns("E:E").ColumnWidth = 24.67

# New header for the inserted column.
$ws.Range("E1").Value = "altname4"

# Democratic Republic of the Congo (row with iso3 "COD") gains a new
# alternate name and a former name. (Former name entered first so the
# shared-string table order matches the source workbook.)
$ws.Range("F49").Value = "Zaire"
$ws.Range("E49").Value = "Congo Democratic Republic"

# Leave the freeze-at-row-1 view untouched (the column insert keeps it),
# and move the active selection the same way the author's session ended up.
$ws.Range("D50").Select() | Out-Null
